$d = $word.ActiveDocument

# 1. Merge "трите имена по документ " + "за самоличност" into one run
#    (identical text, just forces the two adjacent runs to merge into one)
$d.Content.Find.Execute("трите имена по документ за самоличност", $true, $false, $false, $false, $false, $true, 1, $false, "трите имена по документ за самоличност", 2)

# 2. Merge " д" + "ни" + " платен годишен отпуск за " into one run
$d.Content.Find.Execute(" дни платен годишен отпуск за ", $true, $false, $false, $false, $false, $true, 1, $false, " дни платен годишен отпуск за ", 2)

# 3. Prefix the date run with "Дата: "
$d.Content.Find.Execute("date г", $true, $false, $false, $false, $false, $true, 1, $false, "Дата: date г", 2)

# 4. Shrink the leader spacing (106 -> 95 spaces) and merge the trailing
#    "(" / " Подпис )" runs into the same run as the dotted leader.
$old4 = ".                                                                                                          .......................   ( Подпис )"
$new4 = ".                                                                                               .......................   ( Подпис )"
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
